# Generate Report for Archive
#
# The localization status for the sample document moved from "Ready for
# handoff" to "In Translation". That text is shared by the Overview sheet
# (per-language status columns) and each per-language detail sheet
# ("Status" column), so update all of them. The shorter replacement text
# also narrows the "Status"-related columns (their width had been
# auto-fit to the old, longer text).

$wb = $excel.ActiveWorkbook

# --- Update status text: "Ready for handoff" -> "In Translation" ---

# Overview sheet: per-language status shown in columns E (zh-cn) and F (de-de)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

# zh-cn detail sheet: Status column C
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"

# de-de detail sheet: Status column C
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"

# --- Narrow the Status columns to fit the new, shorter text ---
# (auto-fit-style resize after the content shrank)
$wsOverview.Columns.Item(5).ColumnWidth = 12.45
$wsOverview.Columns.Item(6).ColumnWidth = 12.45
$wsZhCn.Columns.Item(3).ColumnWidth = 12.45
$wsDeDe.Columns.Item(3).ColumnWidth = 12.45
